# "Renombre salas éxcel prueba"
# The "Salas" sheet's room names (previously text like "Mezquita-Catedral")
# become plain sequential numbers 1..15. The shared-string table loses the
# 15 now-unused room-name entries, which shifts every other shared-string
# reference in the workbook down accordingly (handled automatically by the
# engine once nothing references those strings any more). We also restore
# two style-3/style-4 row-banding swaps on the Salas sheet, move the active
# sheet/selection from Jueces back to Salas, and tidy the Jueces sheet's
# last column style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Salas")

# Replace each room-name cell (A2:A16) with its plain 1-based row number.
# Writing a numeric literal converts the cell from a shared-string ("s")
# reference to a plain number ("n") and drops the now-unused string from
# sharedStrings.xml automatically.
for ($i = 2; $i -le 16; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# Restore the original alternating row-banding styles (index 3 / index 4)
# for the rows whose banding parity flips in the target file. Copy the
# format from a cell that already carries the desired style so the
# existing style slot is reused instead of allocating a new one.
$bandStyle3 = $ws.Range("A2")   # carries style index 3
$bandStyle4 = $ws.Range("A3")   # carries style index 4

foreach ($r in @(5, 7, 9, 11, 13)) {
    $bandStyle4.Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

$bandStyle3.Copy()
$ws.Cells.Item(12, 1).PasteSpecial(-4122)

# Switch the active sheet/selection back to "Salas" at C11 (it was on
# "Jueces" at C18 before). Jueces' own selection (C18) is left untouched.
$ws.Activate()
$ws.Range("C11").Select()

Write-Output "done"
